$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text columns (Coin name, Link) - no numeric-format concerns
$ws.Range('B41').Value = 'BKEXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('B42').Value = 'KickToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'

# Numeric-looking text columns (Price, Volume%, Hora) - preserve text storage type
$numericTextCells = @(
    @('D2', '254.72'),
    @('E2', '3.52%'),
    @('G2', '11'),
    @('D3', '28.13'),
    @('E3', '-6.72%'),
    @('G3', '11'),
    @('D4', '5.239'),
    @('E4', '1.55%'),
    @('G4', '11'),
    @('D5', '0.05864'),
    @('E5', '1.83%'),
    @('G5', '11'),
    @('D6', '6.713'),
    @('E6', '0.68%'),
    @('G6', '11'),
    @('D7', '0.8645'),
    @('E7', '1.77%'),
    @('G7', '11'),
    @('D8', '0.9810'),
    @('E8', '14.45%'),
    @('G8', '11'),
    @('D9', '0.1410'),
    @('E9', '1.73%'),
    @('G9', '11'),
    @('D10', '0.07179'),
    @('E10', '1.29%'),
    @('G10', '11'),
    @('D11', '0.03177'),
    @('E11', '-2.66%'),
    @('G11', '11'),
    @('D12', '0.09223'),
    @('E12', '-1.57%'),
    @('G12', '11'),
    @('D13', '0.001542'),
    @('E13', '0.83%'),
    @('G13', '11'),
    @('D14', '0.0006071'),
    @('E14', '1.53%'),
    @('G14', '11'),
    @('D15', '0.005792'),
    @('E15', '-1.92%'),
    @('G15', '11'),
    @('D16', '3.511'),
    @('E16', '-0.48%'),
    @('G16', '11'),
    @('E17', '-1.37%'),
    @('G17', '11'),
    @('D18', '2.221'),
    @('E18', '0.16%'),
    @('G18', '11'),
    @('D19', '0.3181'),
    @('E19', '1.86%'),
    @('G19', '11'),
    @('D20', '0.03482'),
    @('E20', '1.96%'),
    @('G20', '11'),
    @('E21', '-1.72%'),
    @('G21', '11'),
    @('D22', '3.544'),
    @('E22', '1.27%'),
    @('G22', '11'),
    @('D23', '0.04167'),
    @('E23', '1.13%'),
    @('G23', '11'),
    @('E24', '-2.08%'),
    @('G24', '11'),
    @('D25', '0.001222'),
    @('E25', '-0.25%'),
    @('G25', '11'),
    @('D26', '0.004796'),
    @('E26', '15.35%'),
    @('G26', '11'),
    @('D27', '0.0001201'),
    @('E27', '0.06%'),
    @('G27', '11'),
    @('D28', '0.0001466'),
    @('E28', '1.20%'),
    @('G28', '11'),
    @('G29', '11'),
    @('G30', '11'),
    @('G31', '11'),
    @('G32', '11'),
    @('G33', '11'),
    @('G34', '11'),
    @('G35', '11'),
    @('G36', '11'),
    @('G37', '11'),
    @('G38', '11'),
    @('G39', '11'),
    @('D40', '0.03810'),
    @('E40', '1.58%'),
    @('G40', '11'),
    @('D41', '0.1101'),
    @('E41', '2.74%'),
    @('G41', '11'),
    @('D42', '0.003823'),
    @('E42', '-32.93%'),
    @('G42', '11'),
    @('D43', '0.002358'),
    @('E43', '-4.11%'),
    @('G43', '11'),
    @('D44', '0.009480'),
    @('E44', '-8.62%'),
    @('G44', '11'),
    @('D45', '0.00005241'),
    @('E45', '-4.31%'),
    @('G45', '11'),
    @('E46', '0.02%'),
    @('G46', '11'),
    @('E47', '69.02%'),
    @('G47', '11'),
    @('D48', '0.002136'),
    @('E48', '-13.38%'),
    @('G48', '11'),
    @('E49', '0.02%'),
    @('G49', '11'),
    @('E50', '0.02%'),
    @('G50', '11'),
    @('G51', '11')
)

foreach ($entry in $numericTextCells) {
    $addr = $entry[0]
    $val = $entry[1]
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}
